$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date serial stored in A5 gets refreshed slightly (same instant,
# re-derived float) as part of this data-retrieval commit.
$ws.Range("A5").Value = 44318.77224265162

# Append the newly retrieved row of job-number stats.
$ws.Range("A6").Value = 44319.77475233001
$ws.Range("B6").Value = 71020
$ws.Range("C6").Value = 59778
$ws.Range("D6").Value = 3279
$ws.Range("E6").Value = 1958
$ws.Range("F6").Value = 1395
$ws.Range("G6").Value = 18546
$ws.Range("H6").Value = 1399
$ws.Range("I6").Value = 798
$ws.Range("J6").Value = 207
